$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" (D) and "is_enabled" (E) columns entirely.
# This shifts the former "order_by" (F) and "rem" (G) columns left into D and E.
$ws.Range("D1:E1").EntireColumn.Delete()
